# Edited CRM accuracy calculation
# The "Batch value" (CRM value measured for this batch) for row 31
# (2021-04-28, Batch #181) was corrected from 2225.47 to 2224.47.
# The "% off" column (D) holds a shared formula that recalculates
# automatically from this input.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRMAccuracyData")

$ws.Range("C31").Value = 2224.4699999999998
